# Update town close year columns
#
# Renames the "2023 ..." / "2024 ..." report headers to generic
# "Prior Year ..." / "Curr. Year ..." labels, moves Card Code / TAXYR /
# TOWNSHIP earlier in the column order, and widens the affected columns
# (which are no longer best-fit because the new header text is longer
# than the old one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row relabel / reorder (row 1, columns A:AA)
# ---------------------------------------------------------------------
$ws.Range("A1").Value  = "PARID"
$ws.Range("B1").Value  = "TAXYR"
$ws.Range("C1").Value  = "TOWNSHIP"
$ws.Range("D1").Value  = "CLASS"
$ws.Range("E1").Value  = "OWN1"
$ws.Range("F1").Value  = "Reason for Change"
$ws.Range("G1").Value  = "WHO"
$ws.Range("H1").Value  = "CARD"
$ws.Range("I1").Value  = "Card Code"
$ws.Range("J1").Value  = "Prior Year Occupancy %"
$ws.Range("K1").Value  = "Curr. Year Occupancy %"
$ws.Range("L1").Value  = "Prior Year % of Curr. Year"
$ws.Range("M1").Value  = "Difference in %"
$ws.Range("N1").Value  = "Prior Year COMDAT VAL"
$ws.Range("O1").Value  = "Curr. Year COMDAT VAL"
$ws.Range("P1").Value  = "Prior Year LMV"
$ws.Range("Q1").Value  = "Prior Year BMV"
$ws.Range("R1").Value  = "Prior Year Total MV"
$ws.Range("S1").Value  = "Prior Year LAV"
$ws.Range("T1").Value  = "Prior Year BAV"
$ws.Range("U1").Value  = "Prior Year Total AV"
$ws.Range("V1").Value  = "Curr. Year LMV"
$ws.Range("W1").Value  = "Curr. Year BMV"
$ws.Range("X1").Value  = "Curr. Year Total MV"
$ws.Range("Y1").Value  = "Curr. Year LAV"
$ws.Range("Z1").Value  = "Curr. Year BAV"
$ws.Range("AA1").Value = "Curr. Year Total AV"

# ---------------------------------------------------------------------
# 2. Column widths
#
# Columns B/C kept their old (pre-reorder) widths; nudge them to the
# post-edit values. Columns J:AA carry the new, much longer labels, so
# they are widened (and are no longer "best fit" - they get an explicit
# width instead). Columns that only shifted by a hundredth of a
# character (font/DPI rendering noise, not a deliberate resize) are left
# alone so their original best-fit autosizing is preserved.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth  = 9.33333333333333   # B  -> 10.109375
$ws.Columns.Item(3).ColumnWidth  = 12.5                # C  -> 13.33203125

$ws.Columns.Item(10).ColumnWidth = 23                  # J  -> 23.88671875
$ws.Columns.Item(11).ColumnWidth = 22.6666666666667    # K  -> 23.5546875
$ws.Columns.Item(12).ColumnWidth = 24.6666666666667    # L  -> 25.5546875
$ws.Columns.Item(13).ColumnWidth = 15.5                # M  -> 16.33203125
$ws.Columns.Item(14).ColumnWidth = 23                  # N  -> 23.77734375
$ws.Columns.Item(15).ColumnWidth = 23.8333333333333    # O  -> 24.6640625
$ws.Columns.Item(16).ColumnWidth = 15.3333333333333    # P  -> 16.21875
$ws.Columns.Item(17).ColumnWidth = 16.3333333333333    # Q  -> 17.21875
$ws.Columns.Item(18).ColumnWidth = 18.8333333333333    # R  -> 19.6640625
$ws.Columns.Item(19).ColumnWidth = 14.5                # S  -> 15.33203125
$ws.Columns.Item(20).ColumnWidth = 14.8333333333333    # T  -> 15.6640625
$ws.Columns.Item(21).ColumnWidth = 18.5                # U  -> 19.33203125
$ws.Columns.Item(22).ColumnWidth = 15                  # V  -> 15.88671875
$ws.Columns.Item(23).ColumnWidth = 15.3333333333333    # W  -> 16.109375
$ws.Columns.Item(24).ColumnWidth = 19.3333333333333    # X  -> 20.109375
$ws.Columns.Item(25).ColumnWidth = 15.8333333333333    # Y  -> 16.6640625
$ws.Columns.Item(26).ColumnWidth = 17.1666666666667    # Z  -> 18
$ws.Columns.Item(27).ColumnWidth = 16.3333333333333    # AA -> 17.21875
